$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.0006408296065709695, 0.002571899574220771, 0.1494219747398047, 0.4942365360607697, 0, 0.6468712399813661)
    3  = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 7.143138311642302)
    4  = @(0.6606524410359556, 0.306821227259698, 22.3905356188092, 10.19245300693656, 0, 33.55046229404141)
    5  = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    6  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    7  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    8  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    9  = @(3.286832544864788, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 0, 3.9712075397307)
    10 = @(1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 3.009163075608874)
    11 = @(0.6606524410359556, 1.655778082260271, 261.3203778131603, 10.19245300693656, 0, 273.8292613433931)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
